$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the new header cells (E1:H1) the same look as the existing ---
# --- bold / centered / wrap-text header cells (reuses style index 1) ---
$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- New header cells (row 1) + new row 2 cells, columns E:H ---
# Order of first-introduction of each distinct new string matters (it drives
# the shared-string table order), so cells are written in the same sequence
# the original author entered them in:
#   email id, firstname, lastname, pujapowar@yopmail.com, role, Super Admin, tester
$ws.Range("E1").Value = "email id"
$ws.Range("F1").Value = "firstname"
$ws.Range("G1").Value = "lastname"
$ws.Range("E2").Value = "pujapowar@yopmail.com"
$ws.Range("H1").Value = "role"
$ws.Range("H2").Value = "Super Admin"
$ws.Range("F2").Value = "tester"
$ws.Range("G2").Value = "lastname"

# --- Update existing row 2 values ---
$ws.Range("A2").Value = "Stage"
$ws.Range("C2").Value = "n8cjtVdQgt"
$ws.Range("D2").Value = "Romania"
$ws.Range("D2").Style = "Normal"

# New hyperlink on E2 (added after the value so the text stays intact)
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:pujapowar@yopmail.com")

# Match the Hyperlink style used for E2
$ws.Range("E2").Style = "Hyperlink"

# --- Row height for header row ---
$ws.Rows.Item(1).RowHeight = 28.8

# --- Selection ---
$ws.Range("F2").Select()
